$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.074.84"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.64%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.750.03"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.16%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.45%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.99"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.05%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.747.85"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.14%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.535"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.04%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.35"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.03%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.463"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.48%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "38.37"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.55%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000249"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +3.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.381.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.40%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.752.76"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.48%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "69.118.87"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.67%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.29"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.98%  "

$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.11"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.99%  "

$ws.Range("E21").Value = "  +18.92%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "493.75"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.39%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000156"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +14.16%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.37"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.33"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.81%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.84%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.42"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.62%  "

$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("E30").Value = "  +7.72%  "

$ws.Range("E31").Value = "  +2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.96"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.12%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "31.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.899.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.46%  "

$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.689.61"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.33%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.24%  "

$ws.Range("E38").Value = "  +2.53%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.26%  "

$ws.Range("E40").Value = "  +2.10%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.81%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.04"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +10.29%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "438.39"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.32%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "48.79"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.23%  "

$ws.Range("E45").Value = "  +2.84%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.47"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "40.55"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "141.75"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.67%  "

$ws.Range("E50").Value = "  +2.52%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.778.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "

